# "RPS Qualifying Source Definitions.xlsx" -- commit: "Fixes RQSD and removes BAU RPS"
$wb = $excel.ActiveWorkbook

$wsBau  = $wb.Worksheets.Item("RQSD-BRQSD")   # BAU RQSD definitions
$wsRqsd = $wb.Worksheets.Item("RQSD-RQSD")    # user-facing RQSD definitions

# --- "Fixes RQSD": a few sources no longer qualify for RPS on the
#     RQSD-RQSD sheet (value 1 -> 0 across the whole B:AE year range) ---
$wsRqsd.Range("B10:AE10").Value = 0   # biomass
$wsRqsd.Range("B24:AE25").Value = 0   # hydrogen combustion turbine, hydrogen combined cycle

# Re-enter the (always-zero) "lignite" formula row on both sheets so Excel
# re-stores it as one shared formula across B14:AE14 (same relative
# formula pattern, =<col>2, that was already there - equivalent to
# selecting B14 and filling right to AE14).
$wsRqsd.Range("B14:AE14").Formula = "=B2"
$wsBau.Range("B14:AE14").Formula = "=B2"

# --- "removes BAU RPS": the BAU sheet is no longer the focused tab; the
#     RQSD-RQSD sheet becomes active/selected instead ---
$wsBau.Activate()
$wsBau.Range("B10").Select()

$wsRqsd.Activate()
$wsRqsd.Range("B24:AE25").Select()
